$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2:C5 keep their cached values but lose the shared formula
$ws.Range("C2").Value = 0.018129770992366411
$ws.Range("C3").Value = 0.007442748091603054
$ws.Range("C4").Value = 0.0093511450381679392
$ws.Range("C5").Value = 0.0091603053435114507

# D2:D5 switch from numeric ratios to the text "U"
$ws.Range("D2").Value = "U"
$ws.Range("D3").Value = "U"
$ws.Range("D4").Value = "U"
$ws.Range("D5").Value = "U"

# New columns E (environ), F (sire), G (dam)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3

# Selection moved to D13
$ws.Range("D13").Select()
